$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update: the existing row 15 (week of 44519) becomes historical
# data that is pushed down to a new row 16, while row 15 is updated in place
# with the new week's figures (44551) and prices.

# 1) Duplicate row 15 into a newly inserted row 16 (keeps all formatting,
#    e.g. the date style on column D).
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(16).Insert()

# 2) Overwrite row 15 with the new week's values.
$ws.Cells.Item(15, 4).Value = 44551   # D15 Fecha

$ws.Cells.Item(15, 13).Value = 400    # M15 Volumen
$ws.Cells.Item(15, 14).Value = 5000   # N15 Precio minimo
$ws.Cells.Item(15, 15).Value = 5500   # O15 Precio maximo
$ws.Cells.Item(15, 16).Value = 5250   # P15 Precio promedio ponderado
$ws.Cells.Item(15, 17).Value = "$/bandeja 12 canastillos 125 gramos"  # Q15 Unidad de comercializacion
$ws.Cells.Item(15, 19).Value = 3500   # S15 Precio $/Kg
$ws.Cells.Item(15, 20).Value = 1.5    # T15 Kg / unidad
